$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.069.09'
$ws.Cells.Item(2, 5).Value = '  +1.58%  '

$ws.Cells.Item(3, 4).Value = '3.427.35'
$ws.Cells.Item(3, 5).Value = '  +1.53%  '

$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.08%  '

$ws.Cells.Item(5, 4).Value = '571.87'
$ws.Cells.Item(5, 5).Value = '  -0.34%  '

$ws.Cells.Item(6, 4).Value = '156.82'
$ws.Cells.Item(6, 5).Value = '  +2.54%  '

$ws.Cells.Item(7, 5).Value = '  -0.10%  '

$ws.Cells.Item(8, 4).Value = '3.423.09'
$ws.Cells.Item(8, 5).Value = '  +1.23%  '

$ws.Cells.Item(9, 4).Value = '0.582'
$ws.Cells.Item(9, 5).Value = '  +10.28%  '

$ws.Cells.Item(10, 4).Value = '7.35'
$ws.Cells.Item(10, 5).Value = '  -0.86%  '

$ws.Cells.Item(11, 4).Value = '0.124'
$ws.Cells.Item(11, 5).Value = '  +4.21%  '

$ws.Cells.Item(12, 4).Value = '0.447'
$ws.Cells.Item(12, 5).Value = '  +2.37%  '

$ws.Cells.Item(13, 4).Value = '4.007.28'
$ws.Cells.Item(13, 5).Value = '  +1.17%  '

$ws.Cells.Item(14, 5).Value = '  -3.07%  '

$ws.Cells.Item(15, 4).Value = '0.0000192'
$ws.Cells.Item(15, 5).Value = '  +6.49%  '

$ws.Cells.Item(16, 4).Value = '27.88'
$ws.Cells.Item(16, 5).Value = '  +3.14%  '

$ws.Cells.Item(17, 4).Value = '64.037.21'
$ws.Cells.Item(17, 5).Value = '  +1.39%  '

$ws.Cells.Item(18, 4).Value = '3.358.51'
$ws.Cells.Item(18, 5).Value = '  -0.17%  '

$ws.Cells.Item(19, 4).Value = '6.45'
$ws.Cells.Item(19, 5).Value = '  +1.75%  '

$ws.Cells.Item(20, 4).Value = '14.32'
$ws.Cells.Item(20, 5).Value = '  +3.00%  '

$ws.Cells.Item(21, 4).Value = '386.47'
$ws.Cells.Item(21, 5).Value = '  +0.21%  '

$ws.Cells.Item(22, 4).Value = '8.16'
$ws.Cells.Item(22, 5).Value = '  -3.00%  '

$ws.Cells.Item(23, 4).Value = '73.48'
$ws.Cells.Item(23, 5).Value = '  +4.21%  '

$ws.Cells.Item(24, 2).Value = 'Polygon'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(24, 4).Value = '0.544'
$ws.Cells.Item(24, 5).Value = '  +1.76%  '

$ws.Cells.Item(25, 2).Value = 'Dai'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(25, 4).Value = '0.997'
$ws.Cells.Item(25, 5).Value = '  -0.43%  '

$ws.Cells.Item(26, 4).Value = '0.0000120'
$ws.Cells.Item(26, 5).Value = '  +23.70%  '

$ws.Cells.Item(27, 4).Value = '9.53'
$ws.Cells.Item(27, 5).Value = '  +2.21%  '

$ws.Cells.Item(28, 4).Value = '0.179'
$ws.Cells.Item(28, 5).Value = '  -0.07%  '

$ws.Cells.Item(29, 5).Value = '  +0.19%  '

$ws.Cells.Item(30, 4).Value = '6.19'
$ws.Cells.Item(30, 5).Value = '  +10.95%  '

$ws.Cells.Item(31, 4).Value = '1.43'
$ws.Cells.Item(31, 5).Value = '  +8.15%  '

$ws.Cells.Item(32, 4).Value = '2.01'
$ws.Cells.Item(32, 5).Value = '  -0.36%  '

$ws.Cells.Item(33, 4).Value = '23.55'
$ws.Cells.Item(33, 5).Value = '  +2.18%  '

$ws.Cells.Item(34, 4).Value = '6.47'
$ws.Cells.Item(34, 5).Value = '  +2.99%  '

$ws.Cells.Item(35, 5).Value = '  +0.15%  '

$ws.Cells.Item(36, 4).Value = '6.97'
$ws.Cells.Item(36, 5).Value = '  +3.89%  '

$ws.Cells.Item(37, 4).Value = '161.11'
$ws.Cells.Item(37, 5).Value = '  +2.21%  '

$ws.Cells.Item(38, 4).Value = '1.47'
$ws.Cells.Item(38, 5).Value = '  -0.52%  '

$ws.Cells.Item(39, 4).Value = '0.0768'
$ws.Cells.Item(39, 5).Value = '  +3.18%  '

$ws.Cells.Item(40, 4).Value = '1.86'
$ws.Cells.Item(40, 5).Value = '  -0.45%  '

$ws.Cells.Item(41, 4).Value = '2.908.89'
$ws.Cells.Item(41, 5).Value = '  +0.58%  '

$ws.Cells.Item(42, 4).Value = '26.99'
$ws.Cells.Item(42, 5).Value = '  -0.96%  '

$ws.Cells.Item(43, 4).Value = '0.0319'
$ws.Cells.Item(43, 5).Value = '  -3.53%  '

$ws.Cells.Item(44, 4).Value = '42.47'
$ws.Cells.Item(44, 5).Value = '  +3.96%  '

$ws.Cells.Item(45, 4).Value = '4.42'
$ws.Cells.Item(45, 5).Value = '  +3.42%  '

$ws.Cells.Item(46, 4).Value = '0.759'
$ws.Cells.Item(46, 5).Value = '  +1.47%  '

$ws.Cells.Item(47, 4).Value = '23.48'
$ws.Cells.Item(47, 5).Value = '  +6.91%  '

$ws.Cells.Item(48, 4).Value = '1.07'
$ws.Cells.Item(48, 5).Value = '  +2.97%  '

$ws.Cells.Item(49, 2).Value = 'dogwifhat'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(49, 4).Value = '2.20'
$ws.Cells.Item(49, 5).Value = '  +21.40%  '

$ws.Cells.Item(50, 2).Value = 'Stellar'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(50, 4).Value = '0.108'
$ws.Cells.Item(50, 5).Value = '  +4.59%  '

$ws.Cells.Item(51, 4).Value = '6.52'
$ws.Cells.Item(51, 5).Value = '  +3.61%  '
